$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new tweets collected by @sharonjcheboi (weeks 2-3 export)
$ws.Range('A278').Value = 'While working to minimise impact of #Covid19 on food security, we need to document what worked, what didn''t, what could have been done different or better or what should have been proporitised more; globally to individual level. #CovidFoodFuture  @sharonjcheboi'
$ws.Range('A280').Value = '#Kenya Justifiably, the government''s focus has shifted to managing #Covid19. Focus should also be drawn towards food security, which is just as important as health security. Poor nutrition promotes susceptibility of an individual to infection. #CovidFoodFuture  @sharonjcheboi'
$ws.Range('A282').Value = '@FAO warned of a massive increase in locusts across East Africa, threatening the already food insecure nations. Control has been slow as borders close due to #Covid19, hence hindering pesticide import.#Covid19FoodFuture @sharonjcheboi'
$ws.Range('A284').Value = '#Kenya #Covid19 expresses the importance of self reliance on food security of a country or household. Encouraging self reliance for those who can access urban agricultural spaces is important. Eg kale, spinach, managu and herbs, potatoes, carrots #Covid19FoodFuture @sharonjcheboi'
$ws.Range('A286').Value = 'To minimise the impact of #Covid19 on food security in #Africa, What actions should be taken to minimise the impact on food distribution? #Covid19FoodFuture @sharonjcheboi'
$ws.Range('A288').Value = 'A crisis such as #Covid19 might be a health crisis. But has the possibility of causing a food crisis. Food systems have been affected and a food crisis within a health crisis should be prevented. Here''s how, by @WorldBank https://bit.ly/2K2DPGV  @sharonjcheboi #Covid19FoodFuture'
$ws.Range('A290').Value = 'This article illustrates the impact of #Covid19 on food joints: Due to stock purchase of products such as tea and coffee, consumers have resorted to home consumption of beverages they would usually consume in hotels. https://bit.ly/2V6c8TP @sharonjcheboi #Covid19FoodFuture'
$ws.Range('A292').Value = 'Farmers globally are facing poor market access as per the poll. Higher value crops which require higher labour are lost before and after harvest, reducing profits. Due to income losses, fewer consumers maybe able to affords these expensive crops. #Covid19FoodFuture  @sharonjcheboi'
$ws.Range('A294').Value = 'Kenya, through @FredMatiangi recently banned uncoordinated food distribution within communities. Was it a good move? Has it provided enough measures to cushion its most vulnerable? Does it disregard food as a basic need? What are your thoughts? #CovidFoodFuture @sharonjcheboi'
$ws.Range('A296').Value = 'In Kenya and all over Africa, many farmers have lost their bulk buyers as schools and markets shut down. Aside from platforms such as @mkulimayoung. How else can farmers cushion themselves and market their produce to minimise complete loss? #CovidFoodFuture @sharonjcheboi'
$ws.Range('A298').Value = 'Yesterday in Kibera, a failed food distribution strategy left people injured and 2 women dead during the chaos. During #covid19, cash transfers to caregivers would be more dignified and adheres to social distancing rules by @WHO #Covid19FoodFuture @sharonjcheboi'
$ws.Range('A300').Value = 'This article illustrates the impact of #Covid19 on food joints. Due to stock purchase of products such as tea and coffee, consumers have resorted to come consumption of beverages they would usually consume in hotels. #Covid19FoodFuture @sharonjcheboi bit.ly/2V6c8TP'
$ws.Range('A302').Value = 'Farmers globally are facing poor market access as per the poll. Higher value crops which require high labour are lost before and after harvest, reducing profits. Due to income losses, fewer consumers may be able to afford these crops. #Covid19FoodFuture @sharonjcheboi'
$ws.Range('A306').Value = 'In Kenya, @wawiranjiru of @Food4Education has set towards providing 6000 households with food staples during #covid19 pandemic. Read more below: https://africasustainabilitymatters.com/social-entrepreneur-28-to-feed-poor-kenyan-homes-amid-coronavirus/amp/?__twitter_impression=true @sharonjcheboi'
$ws.Range('A308').Value = 'In Kenya the @MOH_Kenya announced the first confirmed #covid19 case in Uasin Gishu County. A local farmer has expressed their challenge of accessing labour and inputs due to peoples'' fear of infection and movement restrictions. #covid19foodfuture @sharonjcheboi'
$ws.Range('A310').Value = 'Farm inputs distribution has been restricted and farm input costs have gone higher as a result of higher transportation costs. #Covid19food19future @sharonjcheboi'
$ws.Range('A312').Value = 'Kenyan farmers lost 30 - 35% of 2019 harvest as a result of heavy off season rains. The shortage will very soon be felt in the market. Higher prices for a population that has lost its livelihoods will increase preexistent food insecurity. @sharonjcheboi  #Covid19FoodFuture'
$ws.Range('A314').Value = 'Most food security measures in Kenya focus on production while accessibility is equally important. Due to #COVID19, markets have declined causing massive post-harvest losses. @sharonjcheboi'
$ws.Range('A316').Value = 'Following the heavy rains in January, the recurring locust invasion in February and the #COVID19 pandemic in March, the already vulnerable Kenyan #foodsystem will suffer increasing food insecurity There''s therefore a need to developed resilient systems. @sharonjcheboi'

# Restore the view state: scrolled down with the newly added rows selected
$win = $excel.ActiveWindow
$ws.Range('A278:A318').Select() | Out-Null
$win.ScrollRow = 275
$win.ScrollColumn = 1
